$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 560, pushing existing rows 560:595 down to 561:596
$ws.Rows.Item(560).Insert()

$ws.Cells.Item(560, 1).Value = 6
$ws.Cells.Item(560, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(560, 3).Value = "Metropolitana"
$ws.Cells.Item(560, 4).Value = 45021
$ws.Cells.Item(560, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(560, 5).Value = 13
$ws.Cells.Item(560, 6).Value = 100112043
$ws.Cells.Item(560, 7).Value = "Pepino ensalada"
$ws.Cells.Item(560, 8).Value = "Sin especificar"
$ws.Cells.Item(560, 9).Value = "Primera"
$ws.Cells.Item(560, 10).Value = 570
$ws.Cells.Item(560, 11).Value = 7000
$ws.Cells.Item(560, 12).Value = 8000
$ws.Cells.Item(560, 13).Value = 7596
$ws.Cells.Item(560, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(560, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(560, 16).Value = 127
$ws.Cells.Item(560, 17).Value = 60
$ws.Cells.Item(560, 18).Value = "Hortaliza"
